$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Jumeirah Beach Hotel"
$ws.Range("A3").Value = "London"
$ws.Range("B3").Value = "Grand Plaza Apartments"

$ws.Range("I6").Select()
